# Port input-data updates from eps-eu-ci (through commit 1e366cc) to align modeling
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("About")

# Clean up stale row-level formatting (rows 7, 8 and 9 previously carried a
# row-level style that is no longer used; row 9 becomes fully blank again).
$ws.Rows.Item(7).ClearFormats()
$ws.Rows.Item(8).ClearFormats()
$ws.Rows.Item(9).ClearFormats()

# A12 no longer carries an explicit (redundant) style.
$ws.Range("A12").ClearFormats()

# New "scratch work" that explains the dollar/euro conversion math, added in
# column E next to the existing Large/Medium section headers.
$ws.Range("E20").Value = "want"
$ws.Range("E21").Value = "USD2012/EUR2019"

$c22 = $ws.Range("E22")
$c22.Value = "'="
$c22.ClearFormats()

$ws.Range("E23").Value = "USD2012/USD2019 * USD2019/EUR2019"

$c24 = $ws.Range("E24")
$c24.Value = "'="
$c24.ClearFormats()

$ws.Range("E25").Formula = "=A30*A29"

# Relabel the conversion factors: everything is now expressed in terms of
# USD2012/USD2019/EUR2019 instead of the previous 2012/2019 dollar wording.
$ws.Range("B31").Value = "USD2012 per EUR2019"

$ws.Rows.Item(30).ClearFormats()
$ws.Range("B30").Value = "USD2012 per USD2019"

$ws.Range("A29").ClearFormats()
$ws.Range("B29").Value = "USD2019 per EUR2019"

$ws.Range("G29").Value = "1 USD2019 = 0.8929 EUR2019"
$ws.Range("G30").Value = "1 USD2012 = 1.113 USD2019"

$ws.Range("A29").Formula = "=1/0.8929"
$ws.Range("A30").Formula = "=0.89805"

$ws.Range("F31").Select()

# ---------------------------------------------------------------------------
# Sheet "OCCF-DpSOCU" - formula result changes automatically because it
# references About!A31; only the stale numeric-format style needs clearing.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("OCCF-DpSOCU")
$ws4.Range("B2").ClearFormats()
$ws4.Range("B1").Select()

$wb.Worksheets.Item("About").Activate()
